$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5015
$ws.Range("I62").Value = 5861
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 5861
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -5237
$ws.Range("N62").Value = -4148

$ws.Range("H65").Value = 5015
$ws.Range("I65").Value = 5861
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 29305
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -26185
$ws.Range("N65").Value = -20740

$ws.Range("H76").Value = 50248.43
$ws.Range("I76").Value = 58156.055
$ws.Range("J76").Value = 2802.6667
$ws.Range("K76").Value = 58156.055
$ws.Range("L76").Value = 2802.6667
$ws.Range("M76").Value = -57841.055
$ws.Range("N76").Value = -3432.6667

$ws.Range("H79").Value = 50248.43
$ws.Range("I79").Value = 58156.055
$ws.Range("J79").Value = 2802.6667
$ws.Range("K79").Value = 58156.055
$ws.Range("L79").Value = 2802.6667
$ws.Range("M79").Value = -57064.055
$ws.Range("N79").Value = -4986.6667

$ws.Range("H86").Value = 127426.5
$ws.Range("I86").Value = 201482.4
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 201482.4
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -200359.4
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 127426.5
$ws.Range("I89").Value = 201482.4
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 1007412
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -1001796
$ws.Range("N89").Value = -31232

$ws.Range("H106").Value = 100202000
$ws.Range("I106").Value = 334001.66
$ws.Range("K106").Value = 334001.66
$ws.Range("M106").Value = -333370.66

$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800

$ws.Range("H129").Value = 18710.84
$ws.Range("I129").Value = 550.06665
$ws.Range("K129").Value = 1650.19995
$ws.Range("M129").Value = 3349.80005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2066.8333
$ws.Range("I61").Value = 1150
$ws.Range("J61").Value = 3350.4
$ws.Range("K61").Value = 1150
$ws.Range("L61").Value = 3350.4
$ws.Range("M61").Value = -938
$ws.Range("N61").Value = -3774.4

$ws.Range("H97").Value = 488.80646
$ws.Range("I97").Value = 478.72415
$ws.Range("J97").Value = 635
$ws.Range("K97").Value = 478.72415
$ws.Range("L97").Value = 635
$ws.Range("M97").Value = 17.27584999999999
$ws.Range("N97").Value = -1627

$ws.Range("H136").Value = 2066.8333
$ws.Range("I136").Value = 1150
$ws.Range("J136").Value = 3350.4
$ws.Range("K136").Value = 3450
$ws.Range("L136").Value = 10051.2
$ws.Range("M136").Value = -900
$ws.Range("N136").Value = -15151.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2888.963
$ws.Range("I20").Value = 3273.9473
$ws.Range("J20").Value = 1974.625
$ws.Range("K20").Value = 3273.9473
$ws.Range("L20").Value = 1974.625
$ws.Range("M20").Value = -3026.9473
$ws.Range("N20").Value = -2468.625

$ws.Range("H86").Value = 1759.9375
$ws.Range("I86").Value = 1304.9166
$ws.Range("J86").Value = 3125
$ws.Range("K86").Value = 1304.9166
$ws.Range("L86").Value = 3125
$ws.Range("M86").Value = -181.9166
$ws.Range("N86").Value = -5371

$ws.Range("H89").Value = 1759.9375
$ws.Range("I89").Value = 1304.9166
$ws.Range("J89").Value = 3125
$ws.Range("K89").Value = 6524.583000000001
$ws.Range("L89").Value = 15625
$ws.Range("M89").Value = -908.5830000000005
$ws.Range("N89").Value = -26857

$ws.Range("H94").Value = 972.4231
$ws.Range("I94").Value = 765.61536
$ws.Range("J94").Value = 1179.2307
$ws.Range("K94").Value = 765.61536
$ws.Range("L94").Value = 1179.2307
$ws.Range("M94").Value = -314.61536
$ws.Range("N94").Value = -2081.2307

$ws.Range("H134").Value = 69385.2
$ws.Range("I134").Value = 336004
$ws.Range("J134").Value = 2730.5
$ws.Range("K134").Value = 1008012
$ws.Range("L134").Value = 8191.5
$ws.Range("M134").Value = -1005477
$ws.Range("N134").Value = -13261.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1233.3334
$ws.Range("I16").Value = 850
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -563
$ws.Range("N16").Value = -2574

$ws.Range("H38").Value = 1500
$ws.Range("I38").Value = 1500
$ws.Range("K38").Value = 1500
$ws.Range("M38").Value = -1123

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1289

$ws.Range("H105").Value = 708.6667
$ws.Range("I105").Value = 593.8461
$ws.Range("J105").Value = 895.25
$ws.Range("K105").Value = 593.8461
$ws.Range("L105").Value = 895.25
$ws.Range("M105").Value = 1153.1539
$ws.Range("N105").Value = -4389.25

$ws.Range("H113").Value = 1233.3334
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 850
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1320
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 2823.111
$ws.Range("I132").Value = 2770.4707
$ws.Range("K132").Value = 8311.4121
$ws.Range("M132").Value = -5781.4121

$ws.Range("H134").Value = 1654.4
$ws.Range("I134").Value = 990.2
$ws.Range("K134").Value = 2970.6
$ws.Range("M134").Value = -435.6000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 200.88889
$ws.Range("I14").Value = 200.88889
$ws.Range("K14").Value = 602.6666700000001
$ws.Range("M14").Value = -429.6666700000001

$ws.Range("H36").Value = 636
$ws.Range("I36").Value = 636
$ws.Range("K36").Value = 1908
$ws.Range("M36").Value = -1739

$ws.Range("H75").Value = 33666.668
$ws.Range("I75").Value = 500
$ws.Range("K75").Value = 1500
$ws.Range("M75").Value = -502

$ws.Range("H78").Value = 33666.668
$ws.Range("I78").Value = 500
$ws.Range("K78").Value = 4500
$ws.Range("M78").Value = 492

$ws.Range("H131").Value = 2543.7966
$ws.Range("I131").Value = 17180
$ws.Range("J131").Value = 886.8679
$ws.Range("K131").Value = 51540
$ws.Range("L131").Value = 2660.6037
$ws.Range("M131").Value = -46500
$ws.Range("N131").Value = -12740.6037

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8504066
$ws.Range("I70").Value = 10629097
$ws.Range("J70").Value = 3938.8333
$ws.Range("K70").Value = 10629097
$ws.Range("L70").Value = 3938.8333
$ws.Range("M70").Value = -10628827
$ws.Range("N70").Value = -4478.8333

$ws.Range("H73").Value = 8504066
$ws.Range("I73").Value = 10629097
$ws.Range("J73").Value = 3938.8333
$ws.Range("K73").Value = 10629097
$ws.Range("L73").Value = 3938.8333
$ws.Range("M73").Value = -10628161
$ws.Range("N73").Value = -5810.8333

$ws.Range("H80").Value = 3641.2727
$ws.Range("I80").Value = 2134.5715
$ws.Range("J80").Value = 6278
$ws.Range("K80").Value = 2134.5715
$ws.Range("L80").Value = 6278
$ws.Range("M80").Value = -1136.5715
$ws.Range("N80").Value = -8274

$ws.Range("H83").Value = 3641.2727
$ws.Range("I83").Value = 2134.5715
$ws.Range("J83").Value = 6278
$ws.Range("K83").Value = 10672.8575
$ws.Range("L83").Value = 31390
$ws.Range("M83").Value = -5680.8575
$ws.Range("N83").Value = -41374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 23811130
$ws.Range("I61").Value = 1589.8
$ws.Range("J61").Value = 83334984
$ws.Range("K61").Value = 1589.8
$ws.Range("L61").Value = 83334984
$ws.Range("M61").Value = -1387.8
$ws.Range("N61").Value = -83335388

$ws.Range("H113").Value = 23811130
$ws.Range("I113").Value = 1589.8
$ws.Range("J113").Value = 83334984
$ws.Range("K113").Value = 1589.8
$ws.Range("L113").Value = 83334984
$ws.Range("M113").Value = 580.2
$ws.Range("N113").Value = -83339324

$ws.Range("H122").Value = 2967.923
$ws.Range("I122").Value = 3232.2856
$ws.Range("K122").Value = 9696.856800000001
$ws.Range("M122").Value = -7246.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N109").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0

$ws.Range("H132").Value = 2419.1667
$ws.Range("I132").Value = 1072.3334
$ws.Range("J132").Value = 3766
$ws.Range("K132").Value = 3217.0002
$ws.Range("L132").Value = 11298
$ws.Range("M132").Value = -687.0001999999999
$ws.Range("N132").Value = -16358

Write-Host "Applied all Sheets updates"
